# "added EDDF MP Ground" -- append a new "Ground" strip section below the
# existing "Tower" section (rows 6-8), reusing the same two fill styles
# (orange = row4/"s=1" style, blue = row3/"s=2" style) that the Tower
# section already uses, then resize columns C/D to fit the new content.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section header: "Ground" (mirrors A2 "Tower", no special fill) ---
$ws.Range("A6").Value = "Ground"

# --- Row 7: handover / pushback strip ---
$ws.Range("A7").Value = "READY TO PUSHBACK"
$ws.Range("B7").Value = "PUSHBACK"
$ws.Range("C7").Value = "HANDOVER FROM TOWER"
$ws.Range("D7").Value = "TAXI TO TERMINAL"

# --- Row 8: taxi strip ---
$ws.Range("A8").Value = "TAXI 07L/25R"
$ws.Range("B8").Value = "TAXI 07C/25C"
$ws.Range("C8").Value = "TAXI 07R/25L"
$ws.Range("D8").Value = "TAXI 18"

# Reuse the existing "blue" fill style (same as A3:E3) for A7,B7 and for
# the whole of row 8.
$ws.Range("A3").Copy()
$ws.Range("A7:B7").PasteSpecial(-4122)
$ws.Range("A8:D8").PasteSpecial(-4122)

# Reuse the existing "orange" fill style (same as A4:D4) for C7,D7.
$ws.Range("A4").Copy()
$ws.Range("C7:D7").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Widen column C / narrow column D so the longer "HANDOVER FROM TOWER" /
# shorter "TAXI 18" strings stay best-fit.
$ws.Columns("C").ColumnWidth = 23.1
$ws.Columns("D").ColumnWidth = 16.6

# Move the active selection to D8, matching where the author's cursor
# ended up after typing the new rows.
[void]$ws.Range("D8").Select()
